# Updates the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures on the active worksheet to match the latest scrape, per the
# commit "Updated cryptos list ... with GitHub Actions".
#
# Every one of these cells is stored as plain text (not a number), even
# values such as "1.00" or "0.999" that look numeric. Writing such a
# string straight into .Value would make Excel auto-convert it to a real
# number (dropping the trailing zero, etc.), so for those we prefix the
# value with a leading apostrophe to force text entry and then restore
# the default "Normal" cell style so no stray text-format style sticks
# around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '68.887.89'; ForceText = $False },
    @{ Cell = "E2"; Value = '  +1.50%  '; ForceText = $False },
    @{ Cell = "D3"; Value = '3.762.67'; ForceText = $False },
    @{ Cell = "E3"; Value = '  -0.67%  '; ForceText = $False },
    @{ Cell = "E4"; Value = '  +0.14%  '; ForceText = $False },
    @{ Cell = "D5"; Value = '623.93'; ForceText = $True },
    @{ Cell = "E5"; Value = '  +3.67%  '; ForceText = $False },
    @{ Cell = "D6"; Value = '165.09'; ForceText = $True },
    @{ Cell = "E6"; Value = '  +1.43%  '; ForceText = $False },
    @{ Cell = "D7"; Value = '3.760.54'; ForceText = $False },
    @{ Cell = "E7"; Value = '  -0.67%  '; ForceText = $False },
    @{ Cell = "E8"; Value = '  -0.01%  '; ForceText = $False },
    @{ Cell = "D9"; Value = '0.520'; ForceText = $True },
    @{ Cell = "E9"; Value = '  +1.17%  '; ForceText = $False },
    @{ Cell = "D10"; Value = '0.158'; ForceText = $True },
    @{ Cell = "E10"; Value = '  +0.72%  '; ForceText = $False },
    @{ Cell = "E11"; Value = '  +3.09%  '; ForceText = $False },
    @{ Cell = "D12"; Value = '6.74'; ForceText = $True },
    @{ Cell = "E12"; Value = '  -2.98%  '; ForceText = $False },
    @{ Cell = "D13"; Value = '0.0000245'; ForceText = $True },
    @{ Cell = "E13"; Value = '  +0.02%  '; ForceText = $False },
    @{ Cell = "D14"; Value = '35.53'; ForceText = $True },
    @{ Cell = "E14"; Value = '  +1.55%  '; ForceText = $False },
    @{ Cell = "D15"; Value = '4.406.77'; ForceText = $False },
    @{ Cell = "E15"; Value = '  -0.34%  '; ForceText = $False },
    @{ Cell = "D16"; Value = '3.783.82'; ForceText = $False },
    @{ Cell = "E16"; Value = '  +0.01%  '; ForceText = $False },
    @{ Cell = "D17"; Value = '68.939.93'; ForceText = $False },
    @{ Cell = "E17"; Value = '  +1.58%  '; ForceText = $False },
    @{ Cell = "D18"; Value = '17.62'; ForceText = $True },
    @{ Cell = "E18"; Value = '  -2.99%  '; ForceText = $False },
    @{ Cell = "E19"; Value = '  -1.17%  '; ForceText = $False },
    @{ Cell = "D20"; Value = '7.03'; ForceText = $True },
    @{ Cell = "E20"; Value = '  +0.36%  '; ForceText = $False },
    @{ Cell = "D21"; Value = '466.32'; ForceText = $True },
    @{ Cell = "E21"; Value = '  +1.58%  '; ForceText = $False },
    @{ Cell = "D22"; Value = '9.54'; ForceText = $True },
    @{ Cell = "E22"; Value = '  +1.22%  '; ForceText = $False },
    @{ Cell = "D23"; Value = '0.705'; ForceText = $True },
    @{ Cell = "E23"; Value = '  +2.04%  '; ForceText = $False },
    @{ Cell = "D24"; Value = '0.0000146'; ForceText = $True },
    @{ Cell = "E24"; Value = '  +2.14%  '; ForceText = $False },
    @{ Cell = "D25"; Value = '82.81'; ForceText = $True },
    @{ Cell = "E25"; Value = '  -0.44%  '; ForceText = $False },
    @{ Cell = "D26"; Value = '12.02'; ForceText = $True },
    @{ Cell = "E26"; Value = '  +1.49%  '; ForceText = $False },
    @{ Cell = "D27"; Value = '2.15'; ForceText = $True },
    @{ Cell = "E27"; Value = '  +3.32%  '; ForceText = $False },
    @{ Cell = "D28"; Value = '10.01'; ForceText = $True },
    @{ Cell = "E28"; Value = '  +0.83%  '; ForceText = $False },
    @{ Cell = "D29"; Value = '0.999'; ForceText = $True },
    @{ Cell = "E29"; Value = '  -0.15%  '; ForceText = $False },
    @{ Cell = "D30"; Value = '3.914.06'; ForceText = $False },
    @{ Cell = "E30"; Value = '  -0.57%  '; ForceText = $False },
    @{ Cell = "E31"; Value = '  +2.49%  '; ForceText = $False },
    @{ Cell = "D32"; Value = '2.23'; ForceText = $True },
    @{ Cell = "E32"; Value = '  +2.22%  '; ForceText = $False },
    @{ Cell = "D33"; Value = '7.14'; ForceText = $True },
    @{ Cell = "E33"; Value = '  -0.96%  '; ForceText = $False },
    @{ Cell = "D34"; Value = '28.69'; ForceText = $True },
    @{ Cell = "E34"; Value = '  -0.95%  '; ForceText = $False },
    @{ Cell = "D35"; Value = '0.174'; ForceText = $True },
    @{ Cell = "E35"; Value = '  +19.78%  '; ForceText = $False },
    @{ Cell = "D36"; Value = '1.00'; ForceText = $True },
    @{ Cell = "E36"; Value = '  +0.24%  '; ForceText = $False },
    @{ Cell = "D37"; Value = '3.717.47'; ForceText = $False },
    @{ Cell = "E37"; Value = '  -0.55%  '; ForceText = $False },
    @{ Cell = "D38"; Value = '8.93'; ForceText = $True },
    @{ Cell = "E38"; Value = '  +0.31%  '; ForceText = $False },
    @{ Cell = "E39"; Value = '  +1.92%  '; ForceText = $False },
    @{ Cell = "D40"; Value = '3.34'; ForceText = $True },
    @{ Cell = "E40"; Value = '  +4.85%  '; ForceText = $False },
    @{ Cell = "D41"; Value = '5.81'; ForceText = $True },
    @{ Cell = "E41"; Value = '  +0.35%  '; ForceText = $False },
    @{ Cell = "D42"; Value = '0.967'; ForceText = $True },
    @{ Cell = "E42"; Value = '  -1.22%  '; ForceText = $False },
    @{ Cell = "E43"; Value = '  +0.11%  '; ForceText = $False },
    @{ Cell = "E44"; Value = '  -0.09%  '; ForceText = $False },
    @{ Cell = "D45"; Value = '153.42'; ForceText = $True },
    @{ Cell = "E45"; Value = '  +0.74%  '; ForceText = $False },
    @{ Cell = "D46"; Value = '43.05'; ForceText = $True },
    @{ Cell = "E46"; Value = '  -1.61%  '; ForceText = $False },
    @{ Cell = "D47"; Value = '0.295'; ForceText = $True },
    @{ Cell = "E47"; Value = '  +0.38%  '; ForceText = $False },
    @{ Cell = "D48"; Value = '46.65'; ForceText = $True },
    @{ Cell = "E48"; Value = '  -0.98%  '; ForceText = $False },
    @{ Cell = "D49"; Value = '1.89'; ForceText = $True },
    @{ Cell = "E49"; Value = '  +3.70%  '; ForceText = $False },
    @{ Cell = "D50"; Value = '8.39'; ForceText = $True },
    @{ Cell = "E50"; Value = '  +1.40%  '; ForceText = $False },
    @{ Cell = "D51"; Value = '1.36'; ForceText = $True },
    @{ Cell = "E51"; Value = '  -0.87%  '; ForceText = $False }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to store the value as text
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
